$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.326153333333333
$ws.Range("H2").Value = 3.97846
$ws.Range("I2").Value = 0.01594171638670932
$ws.Range("J2").Value = 0.01594171638670932
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.126464333333334
$ws.Range("N2").Value = 18.379393
$ws.Range("O2").Value = 0.1081098818071741
$ws.Range("P2").Value = 0.1081098818071741
$ws.Range("Q2").Value = 8.124631097197778
$ws.Range("R2").Value = 73.12167987478
$ws.Range("S2").Value = 0.001723457074370635
$ws.Range("T2").Value = 0.001723457074370635

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.326153333333333
$ws.Range("H3").Value = 3.97846
$ws.Range("I3").Value = 0.01594171638670932
$ws.Range("J3").Value = 0.01594171638670932
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 35.73736866666667
$ws.Range("N3").Value = 107.212106
$ws.Range("O3").Value = 0.6306349784216607
$ws.Range("P3").Value = 0.6306349784216608
$ws.Range("Q3").Value = 47.39323058186223
$ws.Range("R3").Value = 426.53907523676
$ws.Range("S3").Value = 0.01005340396953667
$ws.Range("T3").Value = 0.01005340396953667

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.326153333333333
$ws.Range("H4").Value = 3.97846
$ws.Range("I4").Value = 0.01594171638670932
$ws.Range("J4").Value = 0.01594171638670932
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 14.80503233333333
$ws.Range("N4").Value = 44.415097
$ws.Range("O4").Value = 0.2612551397711651
$ws.Range("P4").Value = 0.2612551397711651
$ws.Range("Q4").Value = 19.63374297895778
$ws.Range("R4").Value = 176.70368681062
$ws.Range("S4").Value = 0.004164855342802016
$ws.Range("T4").Value = 0.004164855342802016

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 49.07229233333334
$ws.Range("H5").Value = 147.216877
$ws.Range("I5").Value = 0.5898990314018667
$ws.Range("J5").Value = 0.5898990314018667
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.126464333333334
$ws.Range("N5").Value = 18.379393
$ws.Range("O5").Value = 0.1081098818071741
$ws.Range("P5").Value = 0.1081098818071741
$ws.Range("Q5").Value = 300.6396487350735
$ws.Range("R5").Value = 2705.756838615661
$ws.Range("S5").Value = 0.06377391456302231
$ws.Range("T5").Value = 0.06377391456302231

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 49.07229233333334
$ws.Range("H6").Value = 147.216877
$ws.Range("I6").Value = 0.5898990314018667
$ws.Range("J6").Value = 0.5898990314018667
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 35.73736866666667
$ws.Range("N6").Value = 107.212106
$ws.Range("O6").Value = 0.6306349784216607
$ws.Range("P6").Value = 0.6306349784216608
$ws.Range("Q6").Value = 1753.714602434774
$ws.Range("R6").Value = 15783.43142191296
$ws.Range("S6").Value = 0.3720109629390748
$ws.Range("T6").Value = 0.3720109629390749

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 49.07229233333334
$ws.Range("H7").Value = 147.216877
$ws.Range("I7").Value = 0.5898990314018667
$ws.Range("J7").Value = 0.5898990314018667
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 14.80503233333333
$ws.Range("N7").Value = 44.415097
$ws.Range("O7").Value = 0.2612551397711651
$ws.Range("P7").Value = 0.2612551397711651
$ws.Range("Q7").Value = 726.5168746657856
$ws.Range("R7").Value = 6538.65187199207
$ws.Range("S7").Value = 0.1541141538997696
$ws.Range("T7").Value = 0.1541141538997696

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 32.78916733333333
$ws.Range("H8").Value = 98.367502
$ws.Range("I8").Value = 0.3941592522114239
$ws.Range("J8").Value = 0.3941592522114239
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.126464333333334
$ws.Range("N8").Value = 18.379393
$ws.Range("O8").Value = 0.1081098818071741
$ws.Range("P8").Value = 0.1081098818071741
$ws.Range("Q8").Value = 200.8816641873651
$ws.Range("R8").Value = 1807.934977686286
$ws.Range("S8").Value = 0.04261251016978118
$ws.Range("T8").Value = 0.04261251016978118

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 32.78916733333333
$ws.Range("H9").Value = 98.367502
$ws.Range("I9").Value = 0.3941592522114239
$ws.Range("J9").Value = 0.3941592522114239
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 35.73736866666667
$ws.Range("N9").Value = 107.212106
$ws.Range("O9").Value = 0.6306349784216607
$ws.Range("P9").Value = 0.6306349784216608
$ws.Range("Q9").Value = 1171.798561264357
$ws.Range("R9").Value = 10546.18705137921
$ws.Range("S9").Value = 0.2485706115130493
$ws.Range("T9").Value = 0.2485706115130493

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 32.78916733333333
$ws.Range("H10").Value = 98.367502
$ws.Range("I10").Value = 0.3941592522114239
$ws.Range("J10").Value = 0.3941592522114239
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 14.80503233333333
$ws.Range("N10").Value = 44.415097
$ws.Range("O10").Value = 0.2612551397711651
$ws.Range("P10").Value = 0.2612551397711651
$ws.Range("Q10").Value = 485.4446825530771
$ws.Range("R10").Value = 4369.002142977694
$ws.Range("S10").Value = 0.1029761305285935
$ws.Range("T10").Value = 0.1029761305285935
